$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I4").Value = -0.6759486140615097
$ws.Range("J4").Value = 0.4584227415980057
$ws.Range("K4").Value = 0.4132038225786378
$ws.Range("L4").Value = 2.818383311834202
